# Add an 11th week of larvae-tube labels (codes G101-G110) to the sheet.
#
# Layout: labels run left-to-right, top-to-bottom across columns A, C, E, G
# (4 per row), with columns B, D, F holding a blank "spacer" cell between
# labels. Rows 1-25 already hold codes G1-G100; row 26 exists but is still
# empty in its label columns, rows 27-28 are brand new.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row, col, week, code-number) for every new label, continuing the
# existing T1..T5 week rotation right after G100.
$labels = @(
    @{ Row = 26; Col = 1; Week = 1; Num = 101 },
    @{ Row = 26; Col = 3; Week = 1; Num = 102 },
    @{ Row = 26; Col = 5; Week = 2; Num = 103 },
    @{ Row = 26; Col = 7; Week = 2; Num = 104 },
    @{ Row = 27; Col = 1; Week = 3; Num = 105 },
    @{ Row = 27; Col = 3; Week = 3; Num = 106 },
    @{ Row = 27; Col = 5; Week = 4; Num = 107 },
    @{ Row = 27; Col = 7; Week = 4; Num = 108 },
    @{ Row = 28; Col = 1; Week = 5; Num = 109 },
    @{ Row = 28; Col = 3; Week = 5; Num = 110 }
)

foreach ($label in $labels) {
    $text = "Date :  `nCode : L-T" + $label.Week + "-G" + $label.Num
    $ws.Cells.Item($label.Row, $label.Col).Value = $text
}

# Columns B, D, F are blank spacer cells on every row of the table. Row 26
# already has them; rows 27-28 are new rows and need them created too.
# A bare "" assignment clears/empties a cell instead of materialising it,
# so we type the Excel "treat as text" apostrophe prefix (which resolves
# to the same blank string already used everywhere else) and then reset
# the cell style so it doesn't keep the apostrophe's quote-prefix format.
foreach ($row in 27, 28) {
    foreach ($col in 2, 4, 6) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = "'"
        $cell.Style = "Normal"
    }
}

# Row 28 only needs 10 new codes total, so its last pair (E28/G28) stays
# blank -- but the row's cell grid still gets those two positions touched
# (re-applying the default style materialises an empty cell without
# giving it a value or a new style entry).
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(28, 7).Style = "Normal"

# Re-fit the rows we touched so the multi-line label text we just wrote
# doesn't leave behind an explicit custom row height.
$ws.Range("A26:G28").Rows.AutoFit()
